$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 137
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316295'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Benjamin Bonzi'
$ws.Range("D$r").Value = 'Lorenzo Musetti'
$ws.Range("E$r").Value = 'Gana Lorenzo Musetti'
$ws.Range("F$r").Value = 1.33
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 138
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316266'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Roman Safiullin'
$ws.Range("D$r").Value = 'Holger Rune'
$ws.Range("E$r").Value = 'Gana Holger Rune'
$ws.Range("F$r").Value = 1.33
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 139
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316293'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Stefanos Tsitsipas'
$ws.Range("D$r").Value = 'Fabian Marozsan'
$ws.Range("E$r").Value = 'Gana Stefanos Tsitsipas'
$ws.Range("F$r").Value = 1.91
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 140
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316284'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Terence Atmane'
$ws.Range("D$r").Value = 'Flavio Cobolli'
$ws.Range("E$r").Value = 'Gana Terence Atmane'
$ws.Range("F$r").Value = 3.75
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 141
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316437'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Tatjana Maria'
$ws.Range("D$r").Value = 'Marta Kostyuk'
$ws.Range("E$r").Value = 'Gana Tatjana Maria'
$ws.Range("F$r").Value = 3.75
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 142
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316434'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Victoria Mboko'
$ws.Range("D$r").Value = 'Diana Shnaider'
$ws.Range("E$r").Value = 'Gana Diana Shnaider'
$ws.Range("F$r").Value = 2.1
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 143
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316432'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Maya Joint'
$ws.Range("D$r").Value = 'Beatriz Haddad Maia'
$ws.Range("E$r").Value = 'Gana Beatriz Haddad Maia'
$ws.Range("F$r").Value = 2.38
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 144
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316441'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Clervie Ngounoue'
$ws.Range("D$r").Value = 'Elise Mertens'
$ws.Range("E$r").Value = 'Gana Clervie Ngounoue'
$ws.Range("F$r").Value = 4.33
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 145
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316443'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Olga Danilovic'
$ws.Range("D$r").Value = 'Emma Raducanu'
$ws.Range("E$r").Value = 'Gana Olga Danilovic'
$ws.Range("F$r").Value = 3.2
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 146
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316431'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Anna Kalinskaya'
$ws.Range("D$r").Value = 'Peyton Stearns'
$ws.Range("E$r").Value = 'Gana Peyton Stearns'
$ws.Range("F$r").Value = 2.5
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 147
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316442'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Aryna Sabalenka'
$ws.Range("D$r").Value = 'Marketa Vondrousova'
$ws.Range("E$r").Value = 'Gana Aryna Sabalenka'
$ws.Range("F$r").Value = 1.4
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 148
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14311061'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Carlos Taberner'
$ws.Range("D$r").Value = 'Lukas Neumayer'
$ws.Range("E$r").Value = 'Gana Lukas Neumayer'
$ws.Range("F$r").Value = 3.4
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 149
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14387568'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Oliver Tarvet'
$ws.Range("D$r").Value = 'Henry Searle'
$ws.Range("E$r").Value = 'Gana Henry Searle'
$ws.Range("F$r").Value = 2.75
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 150
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316295'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Benjamin Bonzi'
$ws.Range("D$r").Value = 'Lorenzo Musetti'
$ws.Range("E$r").Value = 'Gana Lorenzo Musetti'
$ws.Range("F$r").Value = 1.33
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 151
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316266'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Roman Safiullin'
$ws.Range("D$r").Value = 'Holger Rune'
$ws.Range("E$r").Value = 'Gana Holger Rune'
$ws.Range("F$r").Value = 1.33
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 152
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316293'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Stefanos Tsitsipas'
$ws.Range("D$r").Value = 'Fabian Marozsan'
$ws.Range("E$r").Value = 'Gana Stefanos Tsitsipas'
$ws.Range("F$r").Value = 1.91
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"

$r = 153
$ws.Range("A$r`:E$r").NumberFormat = "@"
$ws.Range("A$r").Value = '14316284'
$ws.Range("B$r").Value = '2025-08-09'
$ws.Range("C$r").Value = 'Terence Atmane'
$ws.Range("D$r").Value = 'Flavio Cobolli'
$ws.Range("E$r").Value = 'Gana Terence Atmane'
$ws.Range("F$r").Value = 3.75
$ws.Range("G$r`:H$r").NumberFormat = "@"
$ws.Range("G$r`:H$r").Value = "'"
